# EE2 Externalities.pptx -- "Add files via upload" commit
#
# Net logical changes (derived from the canonical OOXML diff):
#   1. Delete the duplicate "Markets are efficient, iff" slide (sldId 277,
#      slide index 12 in the original deck).
#   2. Move the "Resources and Environment" slide (sldId 281) so it sits
#      right before the "Externalities and public goods" section-divider
#      slide (sldId 311) instead of right after it.
#   3. On the "Lindahl Equilibrium" slide, update the embedded photo's
#      alt-text/description from the auto-generated
#      "A person wearing a suit and tie looking at the camera ... " text
#      to "Erik Lindahl".

$p = $ppt.ActivePresentation

# --- 3. Update the picture description on the "Lindahl Equilibrium" slide.
# Do this first, while the original (pre-delete) slide numbering is still
# in effect, so the slide index is unambiguous.
$lindahlSlide = $p.Slides.Item(33)
$photo = $lindahlSlide.Shapes.Item(3)
$photo.AlternativeText = "Erik Lindahl"

# --- 1. Delete the "Markets are efficient, iff" slide at (original) index 12.
$p.Slides.Item(12).Delete()

# --- 2. Move "Resources and Environment" to just before
# "Externalities and public goods". After the deletion above, everything
# from the old index 13 onward shifted down by one, so the
# "Externalities and public goods" divider is now at index 16 and
# "Resources and Environment" is right after it at index 17.
$p.Slides.Item(17).MoveTo(16)
